$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "264.08"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.64"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.204"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06094"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.527"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.726"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8160"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08192"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03365"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03149"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09264"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.921"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001692"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04838"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006239"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006234"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006116"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001101"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001503"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.696"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.274"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3395"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1269"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002685"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04641"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007335"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1123"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003136"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01046"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006093"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7511"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1782"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002103"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01242"
